$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Estadisticos 2P": fill in the second-partial results (previously
# all blank/zeroed placeholders) and add the Promedio column.
# ---------------------------------------------------------------------------
$ws2p = $wb.Worksheets.Item("Estadisticos 2P")

# Row 2 - 4ARHM
$ws2p.Cells.Item(2, 4).Value = 0          # D2 Blancos
$ws2p.Cells.Item(2, 5).Value = 0          # E2 Reprobados
$ws2p.Cells.Item(2, 6).Value = 40         # F2 Aprobados
$ws2p.Cells.Item(2, 7).Value = 100        # G2 Por_Apro
$ws2p.Cells.Item(2, 8).Value = 9.300000000000001  # H2 Promedio

# Row 3 - 6ARHM (DETERMINA LA NOMINA...)
$ws2p.Cells.Item(3, 4).Value = 0          # D3 Blancos
$ws2p.Cells.Item(3, 5).Value = 4          # E3 Reprobados
$ws2p.Cells.Item(3, 6).Value = 35         # F3 Aprobados
$ws2p.Cells.Item(3, 7).Value = 89.73999999999999   # G3 Por_Apro
$ws2p.Cells.Item(3, 8).Value = 8.5        # H3 Promedio

# Row 4 - 6ARHM (INTRODUCCION A LA ECONOMIA)
$ws2p.Cells.Item(4, 4).Value = 0          # D4 Blancos
$ws2p.Cells.Item(4, 5).Value = 2          # E4 Reprobados
$ws2p.Cells.Item(4, 6).Value = 37         # F4 Aprobados
$ws2p.Cells.Item(4, 7).Value = 94.87      # G4 Por_Apro
$ws2p.Cells.Item(4, 8).Value = 8.1        # H4 Promedio

# ---------------------------------------------------------------------------
# Sheet "Estadisticos Final": recompute to reflect the new 2P results.
# ---------------------------------------------------------------------------
$wsFinal = $wb.Worksheets.Item("Estadisticos Final")

# Row 2 - 4ARHM
$wsFinal.Cells.Item(2, 5).Value = 0       # E2 Reprobados
$wsFinal.Cells.Item(2, 6).Value = 40      # F2 Aprobados
$wsFinal.Cells.Item(2, 7).Value = 100     # G2 Por_Apro
$wsFinal.Cells.Item(2, 8).Value = 9.699999999999999  # H2 Promedio

# Row 3 - 6ARHM (DETERMINA LA NOMINA...)
$wsFinal.Cells.Item(3, 8).Value = 8.6     # H3 Promedio

# Row 4 - 6ARHM (INTRODUCCION A LA ECONOMIA)
$wsFinal.Cells.Item(4, 8).Value = 8.199999999999999  # H4 Promedio

# ---------------------------------------------------------------------------
# Sheet "Rescatables": expanded list of students who still need to pass
# (now covering two new students across both subjects, plus the existing
# two students whose "DETERMINA..." entries survive, Kevin Raul's count
# dropping to 2 and Angeles Valeria keeping 1; the INTRODUCCION subject no
# longer carries forward Angeles Valeria/Kevin Raul entries).
# ---------------------------------------------------------------------------
$wsResc = $wb.Worksheets.Item("Rescatables")

$materiaNomina = "DETERMINA LA NÓMINA DEL PERSONAL DE LA ORGANIZACIÓN TOMANDO EN CUENTA LA NORMATIVIDAD LABORAL"
$materiaEconomia = "INTRODUCCIÓN A LA ECONOMÍA"

# New-cell writes below are grouped column-by-column (B, then C, then D)
# across all six rows before moving on, mirroring how the sheet was
# populated originally (pasted columns).

# Column B (Paterno)
$wsResc.Cells.Item(2, 2).Value = "OREA"
$wsResc.Cells.Item(3, 2).Value = "OREA"
$wsResc.Cells.Item(4, 2).Value = "PALMA"
$wsResc.Cells.Item(5, 2).Value = "PALMA"
$wsResc.Cells.Item(6, 2).Value = "RIVERA"
$wsResc.Cells.Item(7, 2).Value = "MORALES"

# Column C (Materno)
$wsResc.Cells.Item(2, 3).Value = "MARTINEZ"
$wsResc.Cells.Item(3, 3).Value = "MARTINEZ"
$wsResc.Cells.Item(4, 3).Value = "RANGEL"
$wsResc.Cells.Item(5, 3).Value = "RANGEL"
$wsResc.Cells.Item(6, 3).Value = "HERNANDEZ"
$wsResc.Cells.Item(7, 3).Value = "ESPARZA"

# Column D (Nombres)
$wsResc.Cells.Item(2, 4).Value = "JOSE MANUEL"
$wsResc.Cells.Item(3, 4).Value = "JOSE MANUEL"
$wsResc.Cells.Item(4, 4).Value = "ROBERTO"
$wsResc.Cells.Item(5, 4).Value = "ROBERTO"
$wsResc.Cells.Item(6, 4).Value = "KEVIN RAUL"
$wsResc.Cells.Item(7, 4).Value = "ANGELES VALERIA"

# Column A (NC)
$wsResc.Cells.Item(2, 1).Value = 21330051920053
$wsResc.Cells.Item(3, 1).Value = 21330051920053
$wsResc.Cells.Item(4, 1).Value = 21330051920242
$wsResc.Cells.Item(5, 1).Value = 21330051920242
$wsResc.Cells.Item(6, 1).Value = 22330051920425
$wsResc.Cells.Item(7, 1).Value = 22330051920218

# Column E (Nombre_Largo / materia)
$wsResc.Cells.Item(2, 5).Value = $materiaNomina
$wsResc.Cells.Item(3, 5).Value = $materiaEconomia
$wsResc.Cells.Item(4, 5).Value = $materiaNomina
$wsResc.Cells.Item(5, 5).Value = $materiaEconomia
$wsResc.Cells.Item(6, 5).Value = $materiaNomina
$wsResc.Cells.Item(7, 5).Value = $materiaNomina

# Column F (Grupo)
$wsResc.Cells.Item(2, 6).Value = "6ARHM"
$wsResc.Cells.Item(3, 6).Value = "6ARHM"
$wsResc.Cells.Item(4, 6).Value = "6ARHM"
$wsResc.Cells.Item(5, 6).Value = "6ARHM"
$wsResc.Cells.Item(6, 6).Value = "6ARHM"
$wsResc.Cells.Item(7, 6).Value = "6ARHM"

# Column G (Reprobadas)
$wsResc.Cells.Item(2, 7).Value = 4
$wsResc.Cells.Item(3, 7).Value = 4
$wsResc.Cells.Item(4, 7).Value = 4
$wsResc.Cells.Item(5, 7).Value = 4
$wsResc.Cells.Item(6, 7).Value = 2
$wsResc.Cells.Item(7, 7).Value = 1
